$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 100000424
$ws.Range("I28").Value = 100000424
$ws.Range("K28").Value = 100000424
$ws.Range("M28").Value = -99999939

$ws.Range("H40").Value = 2466.6667
$ws.Range("J40").Value = 2683.3333
$ws.Range("L40").Value = 2683.3333
$ws.Range("N40").Value = -3033.3333

$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -25968

$ws.Range("H62").Value = 3925
$ws.Range("I62").Value = 1850
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 1850
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -1226
$ws.Range("N62").Value = -7248

$ws.Range("H65").Value = 3925
$ws.Range("I65").Value = 1850
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 9250
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -6130
$ws.Range("N65").Value = -36240

$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 10000
$ws.Range("K69").Value = 30000
$ws.Range("M69").Value = -29126

$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 10000
$ws.Range("K72").Value = 90000
$ws.Range("M72").Value = -85632

$ws.Range("H92").Value = 190
$ws.Range("I92").Value = 190
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 190
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 1058
$ws.Range("N92").ClearContents()

$ws.Range("H98").Value = 1099.6
$ws.Range("I98").Value = 1099.6
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1099.6
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 398.4000000000001
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 1099.6
$ws.Range("I122").Value = 1099.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3298.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -848.7999999999997
$ws.Range("N122").ClearContents()

$ws.Range("H129").Value = 2382.3333
$ws.Range("I129").Value = 1805.125
$ws.Range("K129").Value = 5415.375
$ws.Range("M129").Value = -415.375

$ws.Range("H132").Value = 889.2105
$ws.Range("I132").Value = 905.3333
$ws.Range("K132").Value = 2715.9999
$ws.Range("M132").Value = -185.9998999999998

$ws.Range("H138").Value = 2336.2307
$ws.Range("I138").Value = 1919.5217
$ws.Range("J138").Value = 2935.25
$ws.Range("K138").Value = 5758.5651
$ws.Range("L138").Value = 8805.75
$ws.Range("M138").Value = -618.5650999999998
$ws.Range("N138").Value = -19085.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 111114380
$ws.Range("J2").Value = 4900
$ws.Range("L2").Value = 4900
$ws.Range("N2").Value = -5126

$ws.Range("H32").Value = 5136.1406
$ws.Range("I32").Value = 2606.1226
$ws.Range("K32").Value = 2606.1226
$ws.Range("M32").Value = -2319.1226

$ws.Range("H44").Value = 10995.667
$ws.Range("J44").Value = 10995.667
$ws.Range("L44").Value = 10995.667
$ws.Range("N44").Value = -11971.667

$ws.Range("H63").Value = 1720

$ws.Range("H66").Value = 1720

$ws.Range("H116").Value = 111114380
$ws.Range("J116").Value = 4900
$ws.Range("L116").Value = 4900
$ws.Range("N116").Value = -9488

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 111114380
$ws.Range("J3").Value = 4900
$ws.Range("L3").Value = 4900
$ws.Range("N3").Value = -5128

$ws.Range("H20").Value = 3500
$ws.Range("J20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("N20").Value = -2494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 10142.667
$ws.Range("I15").Value = 11210
$ws.Range("J15").Value = 8008
$ws.Range("K15").Value = 11210
$ws.Range("L15").Value = 8008
$ws.Range("M15").Value = -11040
$ws.Range("N15").Value = -8348

$ws.Range("H31").Value = 2696.6
$ws.Range("I31").Value = 1303.8077
$ws.Range("K31").Value = 1303.8077
$ws.Range("M31").Value = -1008.8077

$ws.Range("H34").Value = 2696.6
$ws.Range("I34").Value = 1303.8077
$ws.Range("K34").Value = 1303.8077
$ws.Range("M34").Value = -1101.8077

$ws.Range("H120").Value = 60326
$ws.Range("J120").Value = 60326
$ws.Range("L120").Value = 60326
$ws.Range("N120").Value = -67584

$ws.Range("H134").Value = 2288.6
$ws.Range("I134").Value = 1785.0667
$ws.Range("K134").Value = 5355.2001
$ws.Range("M134").Value = -2820.2001

$ws.Range("H141").Value = 101918.22
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 101918.22
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 101918.22
$ws.Range("N141").Value = -112278.22
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 995
$ws.Range("I45").Value = 995
$ws.Range("K45").Value = 2985
$ws.Range("M45").Value = -2453

$ws.Range("H132").Value = 2114.551
$ws.Range("I132").Value = 1813.9546
$ws.Range("K132").Value = 16325.5914
$ws.Range("M132").Value = -13795.5914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4000
$ws.Range("I12").Value = 4000
$ws.Range("K12").Value = 4000
$ws.Range("M12").Value = -3860

$ws.Range("H46").Value = 44358.75
$ws.Range("I46").Value = 49129.668
$ws.Range("J46").Value = 30046
$ws.Range("K46").Value = 49129.668
$ws.Range("L46").Value = 30046
$ws.Range("M46").Value = -48973.668
$ws.Range("N46").Value = -30358

$ws.Range("H80").Value = 4137.125
$ws.Range("I80").Value = 3999.5
$ws.Range("J80").Value = 4183
$ws.Range("K80").Value = 3999.5
$ws.Range("L80").Value = 4183
$ws.Range("M80").Value = -3001.5
$ws.Range("N80").Value = -6179

$ws.Range("H83").Value = 4137.125
$ws.Range("I83").Value = 3999.5
$ws.Range("J83").Value = 4183
$ws.Range("K83").Value = 19997.5
$ws.Range("L83").Value = 20915
$ws.Range("M83").Value = -15005.5
$ws.Range("N83").Value = -30899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 898.8
$ws.Range("I55").Value = 624
$ws.Range("J55").Value = 1998
$ws.Range("K55").Value = 624
$ws.Range("L55").Value = 1998
$ws.Range("M55").Value = -451
$ws.Range("N55").Value = -2344

$ws.Range("H132").Value = 4104.9287
$ws.Range("I132").Value = 3905.818
$ws.Range("K132").Value = 11717.454
$ws.Range("M132").Value = -9187.454000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 508.33334
$ws.Range("J2").Value = 508.33334
$ws.Range("L2").Value = 508.33334
$ws.Range("N2").Value = -732.33334

$ws.Range("H18").Value = 12004
$ws.Range("J18").Value = 12004
$ws.Range("L18").Value = 12004
$ws.Range("N18").Value = -12350

$ws.Range("H49").Value = 299666
$ws.Range("I49").Value = 299666
$ws.Range("K49").Value = 299666
$ws.Range("M49").Value = -299436

$ws.Range("H54").Value = 14994
$ws.Range("J54").Value = 14994
$ws.Range("L54").Value = 14994
$ws.Range("N54").Value = -16034

$ws.Range("H132").Value = 2078.2693
$ws.Range("J132").Value = 3089.25
$ws.Range("L132").Value = 9267.75
$ws.Range("N132").Value = -14327.75
